$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "26.271.37"
$ws.Cells.Item(2, 5).Value = "  +0.52%  "

$ws.Cells.Item(3, 4).Value = "1.662.28"
$ws.Cells.Item(3, 5).Value = "  +0.30%  "

$ws.Cells.Item(4, 5).Value = "  +0.77%  "

$ws.Cells.Item(5, 4).Value = "218.46"
$ws.Cells.Item(5, 5).Value = "  +0.15%  "

$ws.Cells.Item(6, 4).Value = "0.5316"
$ws.Cells.Item(6, 5).Value = "  +0.54%  "

$ws.Cells.Item(7, 5).Value = "  +0.73%  "

$ws.Cells.Item(8, 4).Value = "0.2636"
$ws.Cells.Item(8, 5).Value = "  +0.96%  "

$ws.Cells.Item(9, 4).Value = "0.06359"
$ws.Cells.Item(9, 5).Value = "  +0.15%  "

$ws.Cells.Item(10, 4).Value = "20.51"
$ws.Cells.Item(10, 5).Value = "  +0.28%  "

$ws.Cells.Item(11, 4).Value = "0.07830"
$ws.Cells.Item(11, 5).Value = "  +0.49%  "

$ws.Cells.Item(12, 5).Value = "  +1.19%  "

$ws.Cells.Item(13, 4).Value = "1.662.17"
$ws.Cells.Item(13, 5).Value = "  +0.29%  "

$ws.Cells.Item(14, 4).Value = "1.889.84"
$ws.Cells.Item(14, 5).Value = "  +0.27%  "

$ws.Cells.Item(15, 4).Value = "0.5526"
$ws.Cells.Item(15, 5).Value = "  +0.60%  "

$ws.Cells.Item(16, 4).Value = "0.0₅8161"
$ws.Cells.Item(16, 5).Value = "  -0.64%  "

$ws.Cells.Item(17, 4).Value = "65.62"
$ws.Cells.Item(17, 5).Value = "  +0.26%  "

$ws.Cells.Item(19, 4).Value = "4.660"
$ws.Cells.Item(19, 5).Value = "  +1.81%  "

$ws.Cells.Item(20, 4).Value = "193.00"
$ws.Cells.Item(20, 5).Value = "  +0.21%  "

$ws.Cells.Item(21, 4).Value = "10.21"
$ws.Cells.Item(21, 5).Value = "  +1.01%  "

$ws.Cells.Item(22, 5).Value = "  -0.27%  "

$ws.Cells.Item(23, 5).Value = "  +0.72%  "

$ws.Cells.Item(24, 4).Value = "145.13"
$ws.Cells.Item(24, 5).Value = "  +2.29%  "

$ws.Cells.Item(25, 5).Value = "  -2.30%  "

$ws.Cells.Item(26, 4).Value = "7.189"
$ws.Cells.Item(26, 5).Value = "  -1.30%  "

$ws.Cells.Item(27, 4).Value = "16.07"
$ws.Cells.Item(27, 5).Value = "  -0.79%  "

$ws.Cells.Item(28, 4).Value = "1.485"
$ws.Cells.Item(28, 5).Value = "  +2.94%  "

$ws.Cells.Item(29, 5).Value = "  -1.03%  "

$ws.Cells.Item(30, 4).Value = "1.279"
$ws.Cells.Item(30, 5).Value = "  -0.32%  "

$ws.Cells.Item(31, 5).Value = "  +1.50%  "

$ws.Cells.Item(32, 4).Value = "3.270"
$ws.Cells.Item(32, 5).Value = "  +0.07%  "

$ws.Cells.Item(33, 5).Value = "  +1.43%  "

$ws.Cells.Item(34, 5).Value = "  +1.12%  "

$ws.Cells.Item(35, 4).Value = "0.9585"
$ws.Cells.Item(35, 5).Value = "  +0.23%  "

$ws.Cells.Item(36, 4).Value = "2.424"
$ws.Cells.Item(36, 5).Value = "  +0.50%  "

$ws.Cells.Item(37, 4).Value = "0.5779"
$ws.Cells.Item(37, 5).Value = "  +1.07%  "

$ws.Cells.Item(38, 4).Value = "0.01602"
$ws.Cells.Item(38, 5).Value = "  -1.09%  "

$ws.Cells.Item(39, 4).Value = "0.8631"
$ws.Cells.Item(39, 5).Value = "  +1.89%  "

$ws.Cells.Item(40, 4).Value = "5.836"
$ws.Cells.Item(40, 5).Value = "  +0.40%  "

$ws.Cells.Item(41, 5).Value = "  +0.68%  "

$ws.Cells.Item(42, 4).Value = "1.043.08"
$ws.Cells.Item(42, 5).Value = "  +1.70%  "

$ws.Cells.Item(43, 4).Value = "103.90"
$ws.Cells.Item(43, 5).Value = "  +0.72%  "

$ws.Cells.Item(44, 4).Value = "1.802.48"
$ws.Cells.Item(44, 5).Value = "  +0.11%  "

$ws.Cells.Item(45, 4).Value = "57.40"
$ws.Cells.Item(45, 5).Value = "  -0.01%  "

$ws.Cells.Item(46, 4).Value = "0.0₈105"
$ws.Cells.Item(46, 5).Value = "  -5.28%  "

$ws.Cells.Item(47, 4).Value = "1.010"
$ws.Cells.Item(47, 5).Value = "  +0.25%  "

$ws.Cells.Item(48, 4).Value = "0.4383"
$ws.Cells.Item(48, 5).Value = "  +2.08%  "

$ws.Cells.Item(49, 4).Value = "7.961"
$ws.Cells.Item(49, 5).Value = "  +1.86%  "

$ws.Cells.Item(50, 4).Value = "0.05161"
$ws.Cells.Item(50, 5).Value = "  +0.10%  "

$ws.Cells.Item(51, 4).Value = "1.430"
$ws.Cells.Item(51, 5).Value = "  -3.89%  "
